# Actualización desde MV -datos-
# Adds 5 new daily rows (04-10-2021 .. 08-10-2021) to the bottom of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Date (column A), then values for columns C, D, E (column B is not used in this sheet).
$newRows = @(
    @{ Date = "04-10-2021"; C = 4.2;  D = $null; E = 3.39 },
    @{ Date = "05-10-2021"; C = 4.28; D = 3.41;  E = 3.33 },
    @{ Date = "06-10-2021"; C = $null; D = $null; E = 3.31 },
    @{ Date = "07-10-2021"; C = 4.22; D = 3.46;  E = 3.29 },
    @{ Date = "08-10-2021"; C = $null; D = 3.41;  E = 3.27 }
)

$startRow = 190
$r = $startRow
foreach ($item in $newRows) {
    $aCell = $ws.Cells.Item($r, 1)
    $aCell.NumberFormat = "@"
    $aCell.Value = $item.Date
    $aCell.Style = "Normal"

    if ($item.C -ne $null) {
        $ws.Cells.Item($r, 3).Value = $item.C
    }
    if ($item.D -ne $null) {
        $ws.Cells.Item($r, 4).Value = $item.D
    }
    if ($item.E -ne $null) {
        $ws.Cells.Item($r, 5).Value = $item.E
    }

    $r = $r + 1
}
